$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1 = 'bjj compression shorts'
    2 = 'compression running pants'
    3 = 'sit pad'
    4 = 'paintball pants padded'
    5 = 'knee pads breathable'
    6 = 'basketball compression tights for women'
    7 = 'elastic band black mountain'
    8 = 'lacrosse pads youth boys'
    9 = 'position pad'
    10 = 'knee length tights'
    11 = 'capri pants men'
    12 = 'knee pads volleyball black'
    13 = 'softball sliding pants women'
    14 = '6ft basketball'
    15 = 'basketball shorts and pants'
    16 = 'work pants for men construction knee pads'
    17 = 'sliding shorts women softball'
    18 = 'fight shorts wrestling'
    19 = 'youth mesh leggings'
    20 = 'above the knee basketball shorts'
    21 = 'mens yoga leggings'
    22 = 'weight lifting pants for men'
    23 = 'cheap volleyball knee pads'
    24 = 'compression spandex'
    25 = 'yoga pants compression'
    26 = 'mens above the knee shorts'
    27 = 'mens running compression'
    28 = 'cycling pants mens'
    29 = 'knee sleeves basketball men'
    30 = 'softball gear for girls'
    31 = 'cold knee compression'
    32 = 'youth padded compression shorts'
    33 = 'yoga pants for men'
    34 = 'mens spandex tights'
    35 = 'softball protective gear'
    36 = 'soccer sliding shorts'
    37 = 'compression baseball shorts'
    38 = 'long shorts for men below knee'
    39 = 'padded leggings for cycling'
    40 = 'padded volleyball shorts'
    41 = 'hex squat'
    42 = 'youth padded sliding shorts'
    43 = 'knee sleeves bjj'
    44 = 'football pants pads adult'
    45 = 'work pants knee'
    46 = 'cold compression knee'
    47 = '5 pad football girdle'
    48 = 'wrestling sleeve youth'
    49 = 'compression sports pants'
    50 = 'basketball tights for girls'
    51 = 'water pants'
    52 = 'spandex tights men'
    53 = 'boys compression pants black'
    54 = 'hockey tights'
    55 = 'youth hockey compression pants'
    56 = 'men leggings compression'
    57 = 'wrestling kneepads'
    58 = 'kneeling pad gym'
    59 = 'guard shorts'
    60 = 'padded compression shorts men'
    61 = 'softball pants youth'
    62 = 'spandex basketball shorts'
    63 = 'compression shorts men 5 pack'
    64 = 'shorts for men below knee'
    65 = 'mens gym leggings'
    66 = 'compression running leggings'
    67 = 'black mens basketball shorts'
    68 = 'knee pads impact'
    69 = 'paintball pads'
    70 = 'boys compression'
    71 = 'mens volleyball kneepads'
    72 = 'yoga knee pads'
    73 = 'knee work pads'
    74 = 'running capri'
    75 = 'paintball pants for men'
    76 = 'kneepad youth'
    77 = 'polyester capri pants'
    78 = 'man capri pants'
    79 = 'indoor baseball'
    80 = 'softball compression sleeve'
    81 = 'male pads'
    82 = 'high five girls softball pants'
    83 = 'outdoor hockey pants'
    84 = 'basketball floor'
    85 = 'basketball knee sleeve black'
    86 = 'long shorts for men below knee sports'
    87 = 'knee pads for adults'
    88 = 'hockey leggings'
    89 = 'volleyball long knee pads'
    90 = 'lacrosse shorts mens'
    91 = 'mens tights with pouch'
    92 = 'black short baseball pants'
    93 = 'lightweight knee pads'
    94 = 'mens compression pants cold'
    95 = 'knee shorts'
    96 = 'girls sliding pants'
    97 = 'knee pads for work for men'
    98 = 'youth padded leg sleeves for basketball'
    99 = 'gym knee compression'
    100 = 'compression football girdle'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 1).Value = $values[$row]
}
